# "Hoan thanh co ban chuc nang ban dao tao"
# Replace the old "Thong bao tuyen dung" row (row 3) with a new "Thong bao hop
# dinh ky Ban Dao tao" notification, and drop the old "De an cuoc thi NCKH"
# row (row 4) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last notification row (old row 4: "De an cuoc thi NCKH" / BanDaoTao)
$ws.Rows.Item(4).Delete()

# Overwrite row 3 with the new notification's data
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Thông báo họp định kỳ Ban Đào tạo"
$ws.Range("C3").Value = '<p style="text-align: center;"><span style="color: rgb(255, 0, 0);">Thông báo họp định kì ban Đào tạo</span></p>'
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "Ban Đào tạo"
$ws.Range("F3").ClearFormats()
$ws.Range("F3").Value = "11/09/2022 09:15"
$ws.Range("G3").Value = "fdajhfjdashfd"

# Touch the header/footer settings (matches the resulting <headerFooter/> element)
$ws.PageSetup.CenterHeader = ""

# Leave the whole row selected, as after a row delete/edit in the UI
$ws.Rows.Item(3).Select() | Out-Null
